# Payment Reminder module added
# - Rename Sheet1 header row (A1:G1) to new header labels
# - Fill in the CallStatus/CallRemark values for rows 5 and 6
# - Narrow column F to match the new (shorter) "CallRemark" header

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row
$ws.Range("A1").Value = "MemberId"
$ws.Range("B1").Value = "Branch"
$ws.Range("C1").Value = "Fullname"
$ws.Range("D1").Value = "ContactNo"
$ws.Range("E1").Value = "CallStatus"
$ws.Range("F1").Value = "CallRemark"
$ws.Range("G1").Value = "Shift"

# Previously-blank call status / remark cells for rows 5 & 6
$ws.Range("E5").Value = "Called"
$ws.Range("F5").Value = "ok"
$ws.Range("E6").Value = "Called but didn't received"
$ws.Range("F6").Value = "no"

# Column F width shrank from ~17.86 to ~11.86 characters
$ws.Columns.Item(6).ColumnWidth = 11
